$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values repulled from source data.
$updates = @{
    2  = -2
    4  = -1
    5  = -2
    6  = 1
    7  = 4
    8  = -1
    9  = -6
    10 = -5
    11 = 3
    12 = -1
    13 = -1
    14 = -2
    15 = -1
    17 = -4
    18 = 1
    21 = 5
    24 = 0
    26 = -3
    27 = 2
    28 = -1
    29 = -1
    30 = -1
    31 = -1
    32 = -2
    33 = -3
    34 = 1
    35 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
